$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 13036.9
$ws.Range("I34").Value = 13036.9
$ws.Range("K34").Value = 13036.9
$ws.Range("M34").Value = -12833.9

$ws.Range("H36").Value = 13036.9
$ws.Range("I36").Value = 13036.9
$ws.Range("K36").Value = 13036.9
$ws.Range("M36").Value = -12321.9

$ws.Range("H47").Value = 649
$ws.Range("I47").Value = 649
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 649
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 323
$ws.Range("N47").ClearContents()

$ws.Range("H62").Value = 34852.234
$ws.Range("I62").Value = 4766.3335
$ws.Range("J62").Value = 41299.215
$ws.Range("K62").Value = 4766.3335
$ws.Range("L62").Value = 41299.215
$ws.Range("M62").Value = -4142.3335
$ws.Range("N62").Value = -42547.215

$ws.Range("H65").Value = 34852.234
$ws.Range("I65").Value = 4766.3335
$ws.Range("J65").Value = 41299.215
$ws.Range("K65").Value = 23831.6675
$ws.Range("L65").Value = 206496.075
$ws.Range("M65").Value = -20711.6675
$ws.Range("N65").Value = -212736.075

$ws.Range("H98").Value = 4103.206
$ws.Range("I98").Value = 4442.7085
$ws.Range("K98").Value = 4442.7085
$ws.Range("M98").Value = -2944.7085

$ws.Range("H106").Value = 104049.6
$ws.Range("I106").Value = 4499.5557
$ws.Range("K106").Value = 4499.5557
$ws.Range("M106").Value = -3868.5557

$ws.Range("H122").Value = 4103.206
$ws.Range("I122").Value = 4442.7085
$ws.Range("K122").Value = 13328.1255
$ws.Range("M122").Value = -10878.1255

$ws.Range("H137").Value = 2812.3333
$ws.Range("I137").Value = 2445.52
$ws.Range("K137").Value = 7336.559999999999
$ws.Range("M137").Value = -4786.559999999999

$ws.Range("H138").Value = 2468.628
$ws.Range("I138").Value = 2450.8948
$ws.Range("J138").Value = 2482.6667
$ws.Range("K138").Value = 7352.6844
$ws.Range("L138").Value = 7448.000100000001
$ws.Range("M138").Value = -2212.6844
$ws.Range("N138").Value = -17728.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8553.107
$ws.Range("I32").Value = 5255.3022
$ws.Range("K32").Value = 5255.3022
$ws.Range("M32").Value = -4968.3022

$ws.Range("H61").Value = 3857.1428
$ws.Range("I61").Value = 3400
$ws.Range("K61").Value = 3400
$ws.Range("M61").Value = -3188

$ws.Range("H74").Value = 142863140
$ws.Range("I74").Value = 250004640
$ws.Range("J74").Value = 7798
$ws.Range("K74").Value = 250004640
$ws.Range("L74").Value = 7798
$ws.Range("M74").Value = -250003766
$ws.Range("N74").Value = -9546

$ws.Range("H77").Value = 142863140
$ws.Range("I77").Value = 250004640
$ws.Range("J77").Value = 7798
$ws.Range("K77").Value = 1250023200
$ws.Range("L77").Value = 38990
$ws.Range("M77").Value = -1250018832
$ws.Range("N77").Value = -47726

$ws.Range("H132").Value = 1940.6111
$ws.Range("I132").Value = 1234.258
$ws.Range("K132").Value = 3702.774
$ws.Range("M132").Value = -1172.774

$ws.Range("H136").Value = 3857.1428
$ws.Range("I136").Value = 3400
$ws.Range("K136").Value = 10200
$ws.Range("M136").Value = -7650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 410.66666
$ws.Range("I24").Value = 410.66666
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 410.66666
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -175.66666
$ws.Range("N24").ClearContents()

$ws.Range("H134").Value = 2868.524
$ws.Range("I134").Value = 2327.6287
$ws.Range("K134").Value = 6982.886100000001
$ws.Range("M134").Value = -4447.886100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3370.2856
$ws.Range("J8").Value = 4394.8
$ws.Range("L8").Value = 4394.8
$ws.Range("N8").Value = -4674.8

$ws.Range("H22").Value = 1312.375
$ws.Range("I22").Value = 1350
$ws.Range("K22").Value = 1350
$ws.Range("M22").Value = -1000

$ws.Range("H31").Value = 2346.2415
$ws.Range("I31").Value = 1791
$ws.Range("J31").Value = 2864.4666
$ws.Range("K31").Value = 1791
$ws.Range("L31").Value = 2864.4666
$ws.Range("M31").Value = -1496
$ws.Range("N31").Value = -3454.4666

$ws.Range("H34").Value = 2346.2415
$ws.Range("I34").Value = 1791
$ws.Range("J34").Value = 2864.4666
$ws.Range("K34").Value = 1791
$ws.Range("L34").Value = 2864.4666
$ws.Range("M34").Value = -1589
$ws.Range("N34").Value = -3268.4666

$ws.Range("H86").Value = 87327.7
$ws.Range("I86").Value = 209374.75
$ws.Range("K86").Value = 209374.75
$ws.Range("M86").Value = -208251.75

$ws.Range("H89").Value = 87327.7
$ws.Range("I89").Value = 209374.75
$ws.Range("K89").Value = 1046873.75
$ws.Range("M89").Value = -1041257.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13934115
$ws.Range("I4").Value = 24768648
$ws.Range("K4").Value = 74305944
$ws.Range("M4").Value = -74305832

$ws.Range("H14").Value = 240.18182
$ws.Range("I14").Value = 240.18182
$ws.Range("K14").Value = 720.5454599999999
$ws.Range("M14").Value = -547.5454599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7062.857
$ws.Range("I70").Value = 5221.4
$ws.Range("K70").Value = 5221.4
$ws.Range("M70").Value = -4951.4

$ws.Range("H73").Value = 7062.857
$ws.Range("I73").Value = 5221.4
$ws.Range("K73").Value = 5221.4
$ws.Range("M73").Value = -4285.4

$ws.Range("H113").Value = 3693.1
$ws.Range("I113").Value = 2553.5557
$ws.Range("J113").Value = 4625.4546
$ws.Range("K113").Value = 2553.5557
$ws.Range("L113").Value = 4625.4546
$ws.Range("M113").Value = -383.5556999999999
$ws.Range("N113").Value = -8965.454600000001

$ws.Range("H122").Value = 2077.9473
$ws.Range("I122").Value = 1845.1538
$ws.Range("K122").Value = 5535.4614
$ws.Range("M122").Value = -3085.4614

$ws.Range("H126").Value = 9143.951999999999
$ws.Range("I126").Value = 15424.889
$ws.Range("J126").Value = 4433.25
$ws.Range("K126").Value = 46274.667
$ws.Range("L126").Value = 13299.75
$ws.Range("M126").Value = -43804.667
$ws.Range("N126").Value = -18239.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 10999
$ws.Range("I11").Value = 10999
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 10999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -10859
$ws.Range("N11").ClearContents()

$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 500
$ws.Range("J27").Value = 500
$ws.Range("L27").Value = 500
$ws.Range("N27").Value = -714

$ws.Range("H40").Value = 1037305.5
$ws.Range("I40").Value = 1252556
$ws.Range("K40").Value = 1252556
$ws.Range("M40").Value = -1252420

$ws.Range("H61").Value = 10859.777
$ws.Range("I61").Value = 11217.25
$ws.Range("K61").Value = 11217.25
$ws.Range("M61").Value = -11015.25

$ws.Range("H87").Value = 54999
$ws.Range("I87").Value = 54999
$ws.Range("K87").Value = 54999
$ws.Range("M87").Value = -53876

$ws.Range("H90").Value = 54999
$ws.Range("I90").Value = 54999
$ws.Range("K90").Value = 164997
$ws.Range("M90").Value = -159381

$ws.Range("H113").Value = 10859.777
$ws.Range("I113").Value = 11217.25
$ws.Range("K113").Value = 11217.25
$ws.Range("M113").Value = -9047.25

$ws.Range("H132").Value = 45457900
$ws.Range("I132").Value = 47622424
$ws.Range("J132").Value = 2869
$ws.Range("K132").Value = 142867272
$ws.Range("L132").Value = 8607
$ws.Range("M132").Value = -142864742
$ws.Range("N132").Value = -13667

$ws.Range("H133").Value = 298562.5
$ws.Range("J133").Value = 298562.5
$ws.Range("L133").Value = 298562.5
$ws.Range("N133").Value = -303622.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 91499.82000000001
$ws.Range("I17").Value = 111688.22
$ws.Range("K17").Value = 111688.22
$ws.Range("M17").Value = -111516.22

$ws.Range("H81").Value = 15391252
$ws.Range("I81").Value = 2455
$ws.Range("K81").Value = 4910
$ws.Range("M81").Value = -3849

$ws.Range("H84").Value = 15391252
$ws.Range("I84").Value = 2455
$ws.Range("K84").Value = 24550
$ws.Range("M84").Value = -19246

$ws.Range("H113").Value = 3932.225
$ws.Range("J113").Value = 2240
$ws.Range("L113").Value = 6720
$ws.Range("N113").Value = -11060

$ws.Range("H126").Value = 1549
$ws.Range("I126").Value = 1420.6875
$ws.Range("K126").Value = 4262.0625
$ws.Range("M126").Value = -1792.0625
